$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 47: one more solved question (apple and orange on hackerrank)
# Match the date formatting used by the other "Date Solved" cells (no border)
# by copying the format from an existing date cell before writing the value.
$ws.Range("B14").Copy()
$ws.Range("B47").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B47").Value = (Get-Date -Year 2024 -Month 10 -Day 20).Date

$ws.Range("C47").Value = "apple and orange"
$ws.Range("D47").Value = "hackerrank"
$ws.Range("G47").Value = "30days"
$ws.Range("H47").Value = "https://www.hackerrank.com/challenges/apple-and-orange/problem"

$ws.Range("H47").Select()
